$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "22.118.94"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.561.42"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'288.96"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.3799"
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("D8").Value = "'0.3282"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'43.53"
$ws.Range("E9").Value = "  -9.05%  "
$ws.Range("D10").Value = "'1.137"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "'0.07358"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'19.93"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "'5.808"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "'6.897"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "1.557.04"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "'0.00001091"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").Value = "'0.06659"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'85.81"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "'6.448"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'16.12"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").Value = "'11.67"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "22.128.14"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'2.267"
$ws.Range("E25").Value = "  -5.07%  "
$ws.Range("D26").Value = "'2.550"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("D27").Value = "'150.80"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'19.10"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "'4.859"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("D30").Value = "1.731.93"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "'121.37"
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("D32").Value = "'1.120"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "'6.028"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").Value = "'1.811"
$ws.Range("E34").Value = "  -9.39%  "
$ws.Range("D35").Value = "'9.353"
$ws.Range("E35").Value = "  -5.22%  "
$ws.Range("D36").Value = "'0.08168"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").Value = "'5.276"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").Value = "'0.06212"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").Value = "'0.02296"
$ws.Range("E39").Value = "  -6.90%  "
$ws.Range("D40").Value = "'0.2142"
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("D41").Value = "'1.234"
$ws.Range("E41").Value = "  -4.87%  "
$ws.Range("D42").Value = "'11.04"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'0.5978"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("D45").Value = "'13.71"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "'3.745"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "'0.5782"
$ws.Range("E47").Value = "  -5.63%  "
$ws.Range("D48").Value = "'1.973"
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").Value = "'120.48"
$ws.Range("E49").Value = "  -3.87%  "
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("D51").Value = "'0.06980"
$ws.Range("E51").Value = "  -3.40%  "
